$d = $word.ActiveDocument

$d.Content.Find.Execute("643÷2=321, 1", $true, $false, $false, $false, $false, $true, 1, $false, "812÷7=116, 0", 2) | Out-Null
$d.Content.Find.Execute("347÷4=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "160÷7=22, 6", 2) | Out-Null
$d.Content.Find.Execute("222÷5=44, 2", $true, $false, $false, $false, $false, $true, 1, $false, "557÷9=61, 8", 2) | Out-Null
$d.Content.Find.Execute("789÷8=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "594÷8=74, 2", 2) | Out-Null
$d.Content.Find.Execute("539÷3=179, 2", $true, $false, $false, $false, $false, $true, 1, $false, "282÷3=94, 0", 2) | Out-Null
$d.Content.Find.Execute("482÷2=241, 0", $true, $false, $false, $false, $false, $true, 1, $false, "101÷5=20, 1", 2) | Out-Null
$d.Content.Find.Execute("305÷5=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "578÷9=64, 2", 2) | Out-Null
$d.Content.Find.Execute("354÷3=118, 0", $true, $false, $false, $false, $false, $true, 1, $false, "454÷5=90, 4", 2) | Out-Null
$d.Content.Find.Execute("235÷2=117, 1", $true, $false, $false, $false, $false, $true, 1, $false, "818÷3=272, 2", 2) | Out-Null
$d.Content.Find.Execute("300÷2=150, 0", $true, $false, $false, $false, $false, $true, 1, $false, "783÷9=87, 0", 2) | Out-Null
$d.Content.Find.Execute("286÷5=57, 1", $true, $false, $false, $false, $false, $true, 1, $false, "627÷6=104, 3", 2) | Out-Null
$d.Content.Find.Execute("978÷4=244, 2", $true, $false, $false, $false, $false, $true, 1, $false, "494÷5=98, 4", 2) | Out-Null
$d.Content.Find.Execute("269÷8=33, 5", $true, $false, $false, $false, $false, $true, 1, $false, "235÷4=58, 3", 2) | Out-Null
$d.Content.Find.Execute("707÷2=353, 1", $true, $false, $false, $false, $false, $true, 1, $false, "485÷8=60, 5", 2) | Out-Null
$d.Content.Find.Execute("329÷8=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "530÷5=106, 0", 2) | Out-Null
$d.Content.Find.Execute("163÷9=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "524÷6=87, 2", 2) | Out-Null
$d.Content.Find.Execute("783÷6=130, 3", $true, $false, $false, $false, $false, $true, 1, $false, "983÷2=491, 1", 2) | Out-Null
$d.Content.Find.Execute("966÷9=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "954÷7=136, 2", 2) | Out-Null
$d.Content.Find.Execute("689÷3=229, 2", $true, $false, $false, $false, $false, $true, 1, $false, "955÷5=191, 0", 2) | Out-Null
$d.Content.Find.Execute("860÷5=172, 0", $true, $false, $false, $false, $false, $true, 1, $false, "718÷5=143, 3", 2) | Out-Null
$d.Content.Find.Execute("650÷6=108, 2", $true, $false, $false, $false, $false, $true, 1, $false, "794÷2=397, 0", 2) | Out-Null
$d.Content.Find.Execute("747÷8=93, 3", $true, $false, $false, $false, $false, $true, 1, $false, "832÷5=166, 2", 2) | Out-Null
$d.Content.Find.Execute("928÷8=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "249÷4=62, 1", 2) | Out-Null
$d.Content.Find.Execute("490÷5=98, 0", $true, $false, $false, $false, $false, $true, 1, $false, "633÷6=105, 3", 2) | Out-Null
$d.Content.Find.Execute("471÷3=157, 0", $true, $false, $false, $false, $false, $true, 1, $false, "120÷7=17, 1", 2) | Out-Null
